$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update LR-pair (Tgfb3-Tgfbr2) expression/specificity metrics with recomputed TPM values
# Row 2
$ws.Range("G2").Value = 1.103903333333333
$ws.Range("H2").Value = 3.31171
$ws.Range("I2").Value = 0.02393122995918198
$ws.Range("J2").Value = 0.02393122995918198
$ws.Range("M2").Value = 5.575746
$ws.Range("N2").Value = 16.727238
$ws.Range("O2").Value = 0.069238947264747
$ws.Range("P2").Value = 0.069238947264747
$ws.Range("Q2").Value = 6.15508459522
$ws.Range("R2").Value = 55.39576135698
$ws.Range("S2").Value = 0.001656973169124334
$ws.Range("T2").Value = 0.001656973169124334
# Row 3
$ws.Range("G3").Value = 1.103903333333333
$ws.Range("H3").Value = 3.31171
$ws.Range("I3").Value = 0.02393122995918198
$ws.Range("J3").Value = 0.02393122995918198
$ws.Range("O3").Value = 0.8150593598279631
$ws.Range("P3").Value = 0.815059359827963
$ws.Range("Q3").Value = 72.45574215166111
$ws.Range("R3").Value = 652.1016793649501
$ws.Range("S3").Value = 0.01950537297042663
$ws.Range("T3").Value = 0.01950537297042663
# Row 4
$ws.Range("G4").Value = 1.103903333333333
$ws.Range("H4").Value = 3.31171
$ws.Range("I4").Value = 0.02393122995918198
$ws.Range("J4").Value = 0.02393122995918198
$ws.Range("M4").Value = 9.317346333333333
$ws.Range("N4").Value = 27.952039
$ws.Range("O4").Value = 0.11570169290729
$ws.Range("P4").Value = 0.11570169290729
$ws.Range("Q4").Value = 10.28544967518778
$ws.Range("R4").Value = 92.56904707669
$ws.Range("S4").Value = 0.002768883819631011
$ws.Range("T4").Value = 0.002768883819631011
# Row 5
$ws.Range("H5").Value = 70.73212899999999
$ws.Range("I5").Value = 0.5111277390235027
$ws.Range("J5").Value = 0.5111277390235027
$ws.Range("M5").Value = 5.575746
$ws.Range("N5").Value = 16.727238
$ws.Range("O5").Value = 0.069238947264747
$ws.Range("P5").Value = 0.069238947264747
$ws.Range("Q5").Value = 131.461461781078
$ws.Range("R5").Value = 1183.153156029702
$ws.Range("S5").Value = 0.03538994656779767
$ws.Range("T5").Value = 0.03538994656779767
# Row 6
$ws.Range("H6").Value = 70.73212899999999
$ws.Range("I6").Value = 0.5111277390235027
$ws.Range("J6").Value = 0.5111277390235027
$ws.Range("O6").Value = 0.8150593598279631
$ws.Range("P6").Value = 0.815059359827963
$ws.Range("R6").Value = 13927.711093652
$ws.Range("S6").Value = 0.4165994477588103
$ws.Range("T6").Value = 0.4165994477588102
# Row 7
$ws.Range("H7").Value = 70.73212899999999
$ws.Range("I7").Value = 0.5111277390235027
$ws.Range("J7").Value = 0.5111277390235027
$ws.Range("M7").Value = 9.317346333333333
$ws.Range("N7").Value = 27.952039
$ws.Range("O7").Value = 0.11570169290729
$ws.Range("P7").Value = 0.11570169290729
$ws.Range("Q7").Value = 219.6785809290034
$ws.Range("R7").Value = 1977.107228361031
$ws.Range("S7").Value = 0.05913834469689477
$ws.Range("T7").Value = 0.05913834469689477
# Row 8
$ws.Range("G8").Value = 21.446869
$ws.Range("H8").Value = 64.34060699999999
$ws.Range("I8").Value = 0.4649410310173153
$ws.Range("J8").Value = 0.4649410310173154
$ws.Range("M8").Value = 5.575746
$ws.Range("N8").Value = 16.727238
$ws.Range("O8").Value = 0.069238947264747
$ws.Range("P8").Value = 0.069238947264747
$ws.Range("Q8").Value = 119.582294039274
$ws.Range("R8").Value = 1076.240646353466
$ws.Range("S8").Value = 0.03219202752782499
$ws.Range("T8").Value = 0.032192027527825
# Row 9
$ws.Range("G9").Value = 21.446869
$ws.Range("H9").Value = 64.34060699999999
$ws.Range("I9").Value = 0.4649410310173153
$ws.Range("J9").Value = 0.4649410310173154
$ws.Range("O9").Value = 0.8150593598279631
$ws.Range("P9").Value = 0.815059359827963
$ws.Range("Q9").Value = 1407.685585595768
$ws.Range("R9").Value = 12669.17027036191
$ws.Range("S9").Value = 0.3789545390987262
$ws.Range("T9").Value = 0.3789545390987261
# Row 10
$ws.Range("G10").Value = 21.446869
$ws.Range("H10").Value = 64.34060699999999
$ws.Range("I10").Value = 0.4649410310173153
$ws.Range("J10").Value = 0.4649410310173154
$ws.Range("M10").Value = 9.317346333333333
$ws.Range("N10").Value = 27.952039
$ws.Range("O10").Value = 0.11570169290729
$ws.Range("P10").Value = 0.11570169290729
$ws.Range("Q10").Value = 199.8279062386303
$ws.Range("R10").Value = 1798.451156147673
$ws.Range("S10").Value = 0.05379446439076421
$ws.Range("T10").Value = 0.05379446439076421
